$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.183.96'
$ws.Range("E2").Value = '  +0.77%  '

$ws.Range("D3").Value = '2.329.56'
$ws.Range("E3").Value = '  +1.19%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '''310.99'
$ws.Range("E5").Value = '  -1.68%  '

$ws.Range("D6").Value = '''108.15'
$ws.Range("E6").Value = '  +3.87%  '

$ws.Range("E7").Value = '  +0.86%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").Value = '''0.615'
$ws.Range("E9").Value = '  +2.22%  '

$ws.Range("D10").Value = '''40.77'
$ws.Range("E10").Value = '  +3.60%  '

$ws.Range("D11").Value = '''0.0917'
$ws.Range("E11").Value = '  +1.16%  '

$ws.Range("D12").Value = '''8.58'
$ws.Range("E12").Value = '  +1.13%  '

$ws.Range("E13").Value = '  -1.04%  '

$ws.Range("D14").Value = '''1.01'
$ws.Range("E14").Value = '  -0.38%  '

$ws.Range("D15").Value = '''15.47'
$ws.Range("E15").Value = '  +1.02%  '

$ws.Range("D16").Value = '2.681.85'
$ws.Range("E16").Value = '  +1.09%  '

$ws.Range("D17").Value = '2.325.48'
$ws.Range("E17").Value = '  +1.05%  '

$ws.Range("D18").Value = '43.356.78'
$ws.Range("E18").Value = '  +1.30%  '

$ws.Range("E19").Value = '  +0.78%  '

$ws.Range("E20").Value = '  +0.92%  '

$ws.Range("D21").Value = '''13.09'
$ws.Range("E21").Value = '  -7.34%  '

$ws.Range("D22").Value = '''74.11'
$ws.Range("E22").Value = '  +0.11%  '

$ws.Range("D23").Value = '''3.49'
$ws.Range("E23").Value = '  -1.85%  '

$ws.Range("D24").Value = '''267.69'
$ws.Range("E24").Value = '  +1.20%  '

$ws.Range("E25").Value = '  +2.24%  '

$ws.Range("E26").Value = '  -0.12%  '

$ws.Range("D27").Value = '''7.55'
$ws.Range("E27").Value = '  +5.94%  '

$ws.Range("D28").Value = '''11.19'
$ws.Range("E28").Value = '  +2.76%  '

$ws.Range("E29").Value = '  -2.57%  '

$ws.Range("D30").Value = '''38.71'
$ws.Range("E30").Value = '  +3.08%  '

$ws.Range("D31").Value = '''22.61'
$ws.Range("E31").Value = '  +0.84%  '

$ws.Range("D32").Value = '''166.92'
$ws.Range("E32").Value = '  +0.20%  '

$ws.Range("D33").Value = '''0.0890'
$ws.Range("E33").Value = '  +2.00%  '

$ws.Range("D34").Value = '''2.80'
$ws.Range("E34").Value = '  +8.68%  '

$ws.Range("E35").Value = '  +0.70%  '

$ws.Range("D36").Value = '''4.72'
$ws.Range("E36").Value = '  +3.35%  '

$ws.Range("E37").Value = '  -2.41%  '

$ws.Range("D38").Value = '''0.0363'
$ws.Range("E38").Value = '  +3.82%  '

$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '''2.85'
$ws.Range("E39").Value = '  +5.95%  '

$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").Value = '''3.79'
$ws.Range("E40").Value = '  +1.05%  '

$ws.Range("D41").Value = '''1.73'
$ws.Range("E41").Value = '  +9.54%  '

$ws.Range("D42").Value = '''104.56'
$ws.Range("E42").Value = '  +13.52%  '

$ws.Range("D43").Value = '''71.78'
$ws.Range("E43").Value = '  +3.73%  '

$ws.Range("E44").Value = '  +3.18%  '

$ws.Range("D45").Value = '''13.20'
$ws.Range("E45").Value = '  +7.18%  '

$ws.Range("E46").Value = '  +0.07%  '

$ws.Range("D47").Value = '''114.18'
$ws.Range("E47").Value = '  -0.40%  '

$ws.Range("D48").Value = '1.660.41'
$ws.Range("E48").Value = '  -4.37%  '

$ws.Range("D49").Value = '''5.34'
$ws.Range("E49").Value = '  +3.94%  '

$ws.Range("D50").Value = '''8.95'
$ws.Range("E50").Value = '  +1.93%  '

$ws.Range("D51").Value = '''75.89'
$ws.Range("E51").Value = '  -5.71%  '
